# Apply refreshed market-price / leve-profit figures to the Lich Profits workbook.
# Source data is a static snapshot (no formulas in these cells), so each target
# cell is written directly with its new literal value per the scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 481.9091
$ws.Range("J18").Value = 1000.6667
$ws.Range("L18").Value = 1000.6667
$ws.Range("N18").Value = -1568.6667
# Row 76
$ws.Range("H76").Value = 4365.5
$ws.Range("I76").Value = 3650
$ws.Range("J76").Value = 4508.6
$ws.Range("K76").Value = 3650
$ws.Range("L76").Value = 4508.6
$ws.Range("M76").Value = -3335
$ws.Range("N76").Value = -5138.6
# Row 79
$ws.Range("H79").Value = 4365.5
$ws.Range("I79").Value = 3650
$ws.Range("J79").Value = 4508.6
$ws.Range("K79").Value = 3650
$ws.Range("L79").Value = 4508.6
$ws.Range("M79").Value = -2558
$ws.Range("N79").Value = -6692.6
# Row 98
$ws.Range("H98").Value = 2971308.5
$ws.Range("I98").Value = 3638211.2
$ws.Range("J98").Value = 1118800.8
$ws.Range("K98").Value = 3638211.2
$ws.Range("L98").Value = 1118800.8
$ws.Range("M98").Value = -3636713.2
$ws.Range("N98").Value = -1121796.8
# Row 122
$ws.Range("H122").Value = 2971308.5
$ws.Range("I122").Value = 3638211.2
$ws.Range("J122").Value = 1118800.8
$ws.Range("K122").Value = 10914633.6
$ws.Range("L122").Value = 3356402.4
$ws.Range("M122").Value = -10912183.6
$ws.Range("N122").Value = -3361302.4
# Row 132
$ws.Range("H132").Value = 3226.4722
$ws.Range("I132").Value = 1460.3823
$ws.Range("K132").Value = 4381.1469
$ws.Range("M132").Value = -1851.1469
# Row 137
$ws.Range("H137").Value = 1808344.1
$ws.Range("I137").Value = 2976567.2
$ws.Range("J137").Value = 2908.318
$ws.Range("K137").Value = 8929701.600000001
$ws.Range("L137").Value = 8724.954000000002
$ws.Range("M137").Value = -8927151.600000001
$ws.Range("N137").Value = -13824.954
# Row 138
$ws.Range("H138").Value = 1569.75
$ws.Range("I138").Value = 785.85297
$ws.Range("J138").Value = 1973.5758
$ws.Range("K138").Value = 2357.55891
$ws.Range("L138").Value = 5920.7274
$ws.Range("M138").Value = 2782.44109
$ws.Range("N138").Value = -16200.7274
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6463.916
$ws.Range("I32").Value = 5097.1626
$ws.Range("K32").Value = 5097.1626
$ws.Range("M32").Value = -4810.1626
# Row 61
$ws.Range("H61").Value = 4463.7827
$ws.Range("I61").Value = 1981.0588
$ws.Range("J61").Value = 11498.167
$ws.Range("K61").Value = 1981.0588
$ws.Range("L61").Value = 11498.167
$ws.Range("M61").Value = -1769.0588
$ws.Range("N61").Value = -11922.167
# Row 74
$ws.Range("H74").Value = 38348.875
$ws.Range("I74").Value = 42271.2
$ws.Range("K74").Value = 42271.2
$ws.Range("M74").Value = -41397.2
# Row 77
$ws.Range("H77").Value = 38348.875
$ws.Range("I77").Value = 42271.2
$ws.Range("K77").Value = 211356
$ws.Range("M77").Value = -206988
# Row 122
$ws.Range("H122").Value = 5904.107
$ws.Range("I122").Value = 3219.8572
$ws.Range("J122").Value = 13956.857
$ws.Range("K122").Value = 9659.571599999999
$ws.Range("L122").Value = 41870.571
$ws.Range("M122").Value = -7209.571599999999
$ws.Range("N122").Value = -46770.571
# Row 132
$ws.Range("H132").Value = 2260.4746
$ws.Range("I132").Value = 2307.772
$ws.Range("J132").Value = 912.5
$ws.Range("K132").Value = 6923.316
$ws.Range("L132").Value = 2737.5
$ws.Range("M132").Value = -4393.316
$ws.Range("N132").Value = -7797.5
# Row 136
$ws.Range("H136").Value = 4463.7827
$ws.Range("I136").Value = 1981.0588
$ws.Range("J136").Value = 11498.167
$ws.Range("K136").Value = 5943.1764
$ws.Range("L136").Value = 34494.501
$ws.Range("M136").Value = -3393.1764
$ws.Range("N136").Value = -39594.501
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 590.4583
$ws.Range("I107").Value = 545.35297
$ws.Range("K107").Value = 545.35297
$ws.Range("M107").Value = 1374.64703
# Row 134
$ws.Range("H134").Value = 3296.3823
$ws.Range("I134").Value = 2774.2693
$ws.Range("J134").Value = 4993.25
$ws.Range("K134").Value = 8322.8079
$ws.Range("L134").Value = 14979.75
$ws.Range("M134").Value = -5787.8079
$ws.Range("N134").Value = -20049.75
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 91888450
$ws.Range("I31").Value = 1430984.9
$ws.Range("J31").Value = 250189010
$ws.Range("K31").Value = 1430984.9
$ws.Range("L31").Value = 250189010
$ws.Range("M31").Value = -1430689.9
$ws.Range("N31").Value = -250189600
# Row 34
$ws.Range("H34").Value = 91888450
$ws.Range("I34").Value = 1430984.9
$ws.Range("J34").Value = 250189010
$ws.Range("K34").Value = 1430984.9
$ws.Range("L34").Value = 250189010
$ws.Range("M34").Value = -1430782.9
$ws.Range("N34").Value = -250189414
# Row 58
$ws.Range("H58").Value = 2212.7932
$ws.Range("I58").Value = 2215.0386
$ws.Range("K58").Value = 2215.0386
$ws.Range("M58").Value = -2012.0386
# Row 94
$ws.Range("H94").Value = 2420.111
$ws.Range("I94").Value = 1899.5
$ws.Range("J94").Value = 2568.8572
$ws.Range("K94").Value = 1899.5
$ws.Range("L94").Value = 2568.8572
$ws.Range("M94").Value = -1448.5
$ws.Range("N94").Value = -3470.8572
# Row 122
$ws.Range("H122").Value = 969.1539
$ws.Range("I122").Value = 706.125
$ws.Range("J122").Value = 1390
$ws.Range("K122").Value = 2118.375
$ws.Range("L122").Value = 4170
$ws.Range("M122").Value = 331.625
$ws.Range("N122").Value = -9070
# Row 132
$ws.Range("H132").Value = 3289.1738
$ws.Range("I132").Value = 1198.1578
$ws.Range("J132").Value = 13221.5
$ws.Range("K132").Value = 3594.4734
$ws.Range("L132").Value = 39664.5
$ws.Range("M132").Value = -1064.4734
$ws.Range("N132").Value = -44724.5
# Row 136
$ws.Range("H136").Value = 2212.7932
$ws.Range("I136").Value = 2215.0386
$ws.Range("K136").Value = 6645.1158
$ws.Range("M136").Value = -4095.1158
$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 1666.3334
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5189
# Row 72
$ws.Range("H72").Value = 1666.3334
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13944
# Row 107
$ws.Range("H107").Value = 430.2963
$ws.Range("I107").Value = 417.6111
$ws.Range("J107").Value = 455.66666
$ws.Range("K107").Value = 1252.8333
$ws.Range("L107").Value = 1366.99998
$ws.Range("M107").Value = 667.1667
$ws.Range("N107").Value = -5206.999980000001
# Row 122
$ws.Range("H122").Value = 1111.2858
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 7698.8335
$ws.Range("I132").Value = 8738.6
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 78647.40000000001
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -76117.40000000001
$ws.Range("N132").Value = -27560
# Row 140
$ws.Range("H140").Value = 5001678.5
$ws.Range("J140").Value = 1966.5454
$ws.Range("L140").Value = 5899.6362
$ws.Range("N140").Value = -16259.6362
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 24918.475
$ws.Range("I132").Value = 27999.133
$ws.Range("K132").Value = 83997.399
$ws.Range("M132").Value = -81467.399
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 6092.375
$ws.Range("I132").Value = 5249.25
$ws.Range("J132").Value = 6935.5
$ws.Range("K132").Value = 15747.75
$ws.Range("L132").Value = 20806.5
$ws.Range("M132").Value = -13217.75
$ws.Range("N132").Value = -25866.5
# Row 136
$ws.Range("H136").Value = 3191.6274
$ws.Range("I136").Value = 1804.4762
$ws.Range("J136").Value = 9665
$ws.Range("K136").Value = 5413.4286
$ws.Range("L136").Value = 28995
$ws.Range("M136").Value = -2863.4286
$ws.Range("N136").Value = -34095
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 1247.75
$ws.Range("I100").Value = 1247.75
$ws.Range("K100").Value = 2495.5
$ws.Range("M100").Value = -1954.5
# Row 122
$ws.Range("H122").Value = 1973.5
$ws.Range("I122").Value = 1631.3334
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4894.0002
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2444.0002
$ws.Range("N122").Value = -13900
# Row 126
$ws.Range("H126").Value = 1687.3636
$ws.Range("I126").Value = 1682.3334
$ws.Range("J126").Value = 1698.1428
$ws.Range("K126").Value = 5047.0002
$ws.Range("L126").Value = 5094.428400000001
$ws.Range("M126").Value = -2577.0002
$ws.Range("N126").Value = -10034.4284
# Row 132
$ws.Range("H132").Value = 17546132
$ws.Range("I132").Value = 20835312
$ws.Range("J132").Value = 3842.889
$ws.Range("K132").Value = 62505936
$ws.Range("L132").Value = 11528.667
$ws.Range("M132").Value = -62503406
$ws.Range("N132").Value = -16588.667
